# Fixed #366 User content is lost after two generation without edition.
#
# The simple fields `m:usercontent zone1` and `m:endusercontent` were
# serialized as <w:fldSimple w:instr="..."/>. Re-running the generator a
# second time collapsed/duplicated their (empty) cached result, losing the
# user content markers. Expanding them to the explicit begin/instrText/
# separate/end run sequence (as Word normally writes fields with no cached
# result) avoids that problem.

$d = $word.ActiveDocument

function Expand-SimpleField {
    param(
        [int]$FieldIndex,
        [string]$InstrText
    )

    $field = $d.Fields.Item($FieldIndex)

    # Remember where the field starts (its containing paragraph/run
    # position) *before* deleting it, so we can re-insert the expanded
    # field markers in exactly the same spot.
    $insertAt = $field.Code.Start - 1

    # Drop the collapsed <w:fldSimple> field entirely.
    $field.Delete()

    # Escape for safe embedding inside the XML instruction text.
    $escaped = $InstrText -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

    $ooxml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
        '<w:r><w:instrText>' + $escaped + '</w:instrText></w:r>' +
        '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
        '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'

    # Insert into a *collapsed* range at the original field position so the
    # new runs are spliced into the existing paragraph instead of replacing
    # it outright.
    $target = $d.Range($insertAt, $insertAt)
    $target.InsertXML($ooxml)
}

# Both fields collapse back to Fields.Item(1) once the previous one has
# been expanded (expanding doesn't add a "real" field, so the remaining
# simple field is always first in the collection).
Expand-SimpleField -FieldIndex 1 -InstrText "m:usercontent zone1"
Expand-SimpleField -FieldIndex 1 -InstrText "m:endusercontent"
